# Updated cryptos list on Fri Jun 14 09:37:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text in the source data
# (e.g. "67.090.31", "  -0.45%  "). Excel's COM layer auto-coerces numeric
# looking strings into real numbers, which would corrupt values such as
# "1.81" or drop the thousands separators in "67.038.76". Force the
# cells to Text format before writing so the literal string is preserved.
$textCells = @(
  "D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","E8",
  "D9","E9","E10","D11","E11","E12","E13","D14","D15","E15","D16","E16",
  "D17","E17","E18","D19","E19","D20","E20","D21","E21","D22","E22",
  "D23","E23","D24","E24","D25","E25","E26","E27","D28","E28","E29",
  "E30","E31","E32","E33","D34","E34","D35","E35",
  "D36","E36","D37","E37","D38","E38","E39","E40",
  "D41","E41","D42","E42","D43","E43","E44","E45","D46","E46","D47","E47",
  "E48","E49","D50","E50","D51","E51"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.038.76"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.519.93"
$ws.Range("E3").Value = "  +1.13%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "608.69"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6 - Solana
$ws.Range("D6").Value = "148.19"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.517.23"
$ws.Range("E7").Value = "  +1.17%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -1.35%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.28%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "7.93"
$ws.Range("E11").Value = "  +5.16%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.56%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.82%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "32.04"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.118.29"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.509.60"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.196.16"
$ws.Range("E17").Value = "  -0.60%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.28%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "10.73"
$ws.Range("E19").Value = "  +9.59%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  -0.57%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "15.33"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "438.31"
$ws.Range("E22").Value = "  -1.60%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  -2.41%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "79.42"
$ws.Range("E24").Value = "  +1.64%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.665.08"
$ws.Range("E25").Value = "  +1.01%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -3.38%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  -1.66%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  -4.44%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.62%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -3.20%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -2.07%  "

# Row 33 - Binance-PegBSC-USD
$ws.Range("E33").Value = "  +0.04%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "25.48"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35 - RenzoRestakedETH
$ws.Range("D35").Value = "3.514.05"
$ws.Range("E35").Value = "  +1.04%  "

# Row 36 - now NEARProtocol (was ImmutableX)
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  -2.55%  "

# Row 37 - now ImmutableX (was NEARProtocol)
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  -2.20%  "

# Row 38 - Aptos
$ws.Range("D38").Value = "8.04"
$ws.Range("E38").Value = "  +1.14%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.02%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.05%  "

# Row 41 - now Hedera (was Monero)
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.0895"
$ws.Range("E41").Value = "  +0.35%  "

# Row 42 - now Monero (was Hedera)
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "172.65"
$ws.Range("E42").Value = "  -2.37%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "5.44"
$ws.Range("E43").Value = "  +0.61%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -9.27%  "

# Row 45 - Mantle
$ws.Range("E45").Value = "  +0.83%  "

# Row 46 - OKB
$ws.Range("D46").Value = "46.02"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "28.17"
$ws.Range("E47").Value = "  -6.26%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  -0.51%  "

# Row 49 - dogwifhat
$ws.Range("E49").Value = "  -2.64%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "7.48"
$ws.Range("E50").Value = "  -1.31%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "0.990"
$ws.Range("E51").Value = "  +0.68%  "
